# The document's header/footer both carry the Pearson and BTec logos as
# inline pictures. The authoring tool re-labelled the auto-generated
# "imageN.ext" display names on these pictures (no visual/content change):
#   - the Pearson logo (appears in both the primary and first-page footers)
#     goes from "image1.png" to "image2.png"
#   - the BTec logo (first-page header) goes from "image2.jpg" to "image1.jpg"

$d = $word.ActiveDocument
$section = $d.Sections(1)

$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

# BTec logo lives in the first-page header.
$btecHeader = $section.Headers($wdHeaderFooterFirstPage)
$btecLogo = $btecHeader.Range.InlineShapes(1)
$btecLogo.Name = "image1.jpg"

# Pearson logo lives in both the primary and the first-page footers.
$primaryFooter = $section.Footers($wdHeaderFooterPrimary)
$pearsonLogo1 = $primaryFooter.Range.InlineShapes(1)
$pearsonLogo1.Name = "image2.png"

$firstPageFooter = $section.Footers($wdHeaderFooterFirstPage)
$pearsonLogo2 = $firstPageFooter.Range.InlineShapes(1)
$pearsonLogo2.Name = "image2.png"
